$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1: update the "Last status check" timestamp
$ws.Range("F1").Value = "Last status check on: 16.02.2022 12:00"

# Row 3 (Tesco): price moved from 35.9 to 36.5 (a +0.6 increase)
$ws.Range("B3").Value = 36.5
$ws.Range("C3").Value = 35.9

# D3 used to hold a numeric delta; it is now a text "+0.6" label.
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "+0.6"
$ws.Range("D3").ClearFormats()

# E3 used to hold a numeric date serial (with a date-time style); it is now
# a plain text timestamp with the default (no) style.
$ws.Range("E3").ClearFormats()
$ws.Range("E3").Value = "2022-02-16 12:00:18"
